$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in columns A, Q, R between row 4 and row 6
$a4 = $ws.Range("A4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

$a6 = $ws.Range("A6").Value2
$q6 = $ws.Range("Q6").Value2
$r6 = $ws.Range("R6").Value2

$ws.Range("A4").Value2 = $a6
$ws.Range("Q4").Value2 = $q6
$ws.Range("R4").Value2 = $r6

$ws.Range("A6").Value2 = $a4
$ws.Range("Q6").Value2 = $q4
$ws.Range("R6").Value2 = $r4

# Swap values in columns A, Q, R between row 5 and row 7
$a5 = $ws.Range("A5").Value2
$q5 = $ws.Range("Q5").Value2
$r5 = $ws.Range("R5").Value2

$a7 = $ws.Range("A7").Value2
$q7 = $ws.Range("Q7").Value2
$r7 = $ws.Range("R7").Value2

$ws.Range("A5").Value2 = $a7
$ws.Range("Q5").Value2 = $q7
$ws.Range("R5").Value2 = $r7

$ws.Range("A7").Value2 = $a5
$ws.Range("Q7").Value2 = $q5
$ws.Range("R7").Value2 = $r5
